# Add/update metadata report for Akurana
# Appends a new data row (row 3) to the "Metadata Report" sheet, extending
# the sheet's dimension through column AO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a real number (year); the rest are text.
$ws.Cells.Item(3, 1).Value = 2024

$textValues = @{
    2  = "DEC"
    3  = "31/12-01/12"
    4  = "-"
    5  = "Akurana"
    6  = "-"
    7  = "-"
    8  = "-"
    9  = "-"
    10 = "-"
    11 = "-"
    12 = "-"
    13 = "-"
    14 = "-"
    15 = "-"
    16 = "-"
    17 = "-"
    18 = "-"
    19 = "-"
    20 = "-"
    21 = "-"
    22 = "-"
    23 = "-"
    24 = "-"
    25 = "-"
    26 = "-"
    27 = "-"
    28 = "-"
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
    35 = "-"
    36 = "-"
    37 = "-"
    38 = "-"
    39 = "-"
    40 = "-"
    41 = "-"
}

foreach ($col in $textValues.Keys) {
    $ws.Cells.Item(3, $col).Value = $textValues[$col]
}
